# Append new scraped rows (2026-01-08 12:41:00 fetch) into the "ランサーズ" sheet.
# The existing 5 rows (priority scores 75/33/30/18/13) are kept, re-sorted by
# score together with 3 brand-new items, and every row is re-stamped with the
# new fetch timestamp. This mirrors the source tool's behaviour: it rewrites
# the data rows top to bottom and only *adds* new Hyperlink relationships for
# the rows that land beyond the sheet's previous extent (row 6) - it does not
# rebind the 5 hyperlinks that already existed, even though the text shown in
# those cells moves to a different row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$stamp = "2026-01-08 12:41:00"

# --- Row 2 : unchanged item (score 75), just re-stamped -------------------
$ws.Cells.Item(2,1).Value = $stamp
$ws.Cells.Item(2,2).Value = "【法人歓迎】プローバステージ制御ソフト開発の見積依頼"
$ws.Cells.Item(2,3).Value = "システム開発"
$ws.Cells.Item(2,4).Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Cells.Item(2,5).Value = "期限情報なし"
$ws.Cells.Item(2,6).Value = "https://www.lancers.jp/work/detail/5467295"
$ws.Cells.Item(2,7).Value = 75
$ws.Cells.Item(2,8).Value = "◆開発"

# --- Row 3 : NEW item (score 38) -------------------------------------------
$ws.Cells.Item(3,1).Value = $stamp
$ws.Cells.Item(3,2).Value = "イベントサイトのWeb制作(決済機能付き)依頼"
$ws.Cells.Item(3,3).Value = "システム開発"
$ws.Cells.Item(3,4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(3,5).Value = "期限情報なし"
$ws.Cells.Item(3,6).Value = "https://www.lancers.jp/work/detail/5467460"
$ws.Cells.Item(3,7).Value = 38
$ws.Cells.Item(3,8).Value = "◇サイト"

# --- Row 4 : previously row 3 (score 33) -----------------------------------
$ws.Cells.Item(4,1).Value = $stamp
$ws.Cells.Item(4,2).Value = "【急募】社内Webアプリの修正・再構築依頼"
$ws.Cells.Item(4,3).Value = "システム開発"
$ws.Cells.Item(4,4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(4,5).Value = "期限情報なし"
$ws.Cells.Item(4,6).Value = "https://www.lancers.jp/work/detail/5467384"
$ws.Cells.Item(4,7).Value = 33
$ws.Cells.Item(4,8).Value = "◇アプリ"

# --- Row 5 : NEW item (score 30) -------------------------------------------
$ws.Cells.Item(5,1).Value = $stamp
$ws.Cells.Item(5,2).Value = "iPhoneのブラウザ要素の書き換えアプリ作成"
$ws.Cells.Item(5,3).Value = "システム開発"
$ws.Cells.Item(5,4).Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Cells.Item(5,5).Value = "期限情報なし"
$ws.Cells.Item(5,6).Value = "https://www.lancers.jp/work/detail/5467578"
$ws.Cells.Item(5,7).Value = 30
$ws.Cells.Item(5,8).Value = "◇アプリ"

# --- Row 6 : previously row 4 (score 30) -----------------------------------
$ws.Cells.Item(6,1).Value = $stamp
$ws.Cells.Item(6,2).Value = "進行管理およびチームディレクションを担当"
$ws.Cells.Item(6,3).Value = "システム開発"
$ws.Cells.Item(6,4).Value = "~ 5,000 円 / 固定"
$ws.Cells.Item(6,5).Value = "期限情報なし"
$ws.Cells.Item(6,6).Value = "https://www.lancers.jp/work/detail/5418064"
$ws.Cells.Item(6,7).Value = 30
$ws.Cells.Item(6,8).Value = "◇管理"

# --- Row 7 : previously row 5 (score 18), beyond the old extent -----------
$ws.Cells.Item(7,1).Value = $stamp
$ws.Cells.Item(7,2).Value = "【急募】cloudflare導入の経験者を探しています!"
$ws.Cells.Item(7,3).Value = "システム開発"
$ws.Cells.Item(7,4).Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item(7,5).Value = "期限情報なし"
$ws.Cells.Item(7,6).Value = "https://www.lancers.jp/work/detail/5467334"
$ws.Cells.Item(7,7).Value = 18

# --- Row 8 : previously row 6 (score 13), beyond the old extent -----------
$ws.Cells.Item(8,1).Value = $stamp
$ws.Cells.Item(8,2).Value = "電気点火装置の回路図作成依頼"
$ws.Cells.Item(8,3).Value = "システム開発"
$ws.Cells.Item(8,4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(8,5).Value = "期限情報なし"
$ws.Cells.Item(8,6).Value = "https://www.lancers.jp/work/detail/5466994"
$ws.Cells.Item(8,7).Value = 13

# --- Row 9 : NEW item (score 10), beyond the old extent --------------------
$ws.Cells.Item(9,1).Value = $stamp
$ws.Cells.Item(9,2).Value = "ドメインの移行をして欲しい"
$ws.Cells.Item(9,3).Value = "システム開発"
$ws.Cells.Item(9,4).Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Cells.Item(9,5).Value = "期限情報なし"
$ws.Cells.Item(9,6).Value = "https://www.lancers.jp/work/detail/5467598"
$ws.Cells.Item(9,7).Value = 10

# --- Hyperlinks: F2:F6 already carry live hyperlinks from before this edit
# (their text just moved underneath them, untouched by design - matches the
# source tool's behaviour). Only F7:F9 are brand new cells and need brand
# new Hyperlink relationships, styled like the rest of the URL column.
$ws.Hyperlinks.Add($ws.Cells.Item(7,6), "https://www.lancers.jp/work/detail/5467334")
$ws.Cells.Item(7,6).Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Cells.Item(8,6), "https://www.lancers.jp/work/detail/5466994")
$ws.Cells.Item(8,6).Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Cells.Item(9,6), "https://www.lancers.jp/work/detail/5467598")
$ws.Cells.Item(9,6).Style = "Hyperlink"
